# Added 0 as an alternative option for barcode_offset
$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force the cell to hold `text` as a real text value (not a number),
    # the same way a leading apostrophe does in the Excel UI, then strip
    # the resulting "text" number-format override so no stray style is
    # left on the cell (matches plain, style-less shared-string cells).
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# 1. Prepend "0" as a new first entry to the barcode_offset lookup sheet,
#    pushing the existing 5 entries down by one row. Final list becomes:
#    0, 0,38,76, Not applicable, 1,27, 8, 10,48,86
$boSheet = $wb.Worksheets.Item("barcode_offset")

$existing = @()
for ($r = 1; $r -le 5; $r++) {
    $existing += ,$boSheet.Cells.Item($r, 1).Text
}

# Shift the existing 5 values down into rows 2-6 (bottom-up so nothing is
# overwritten before it has been captured).
for ($r = 5; $r -ge 1; $r--) {
    Set-TextValue $boSheet.Cells.Item($r + 1, 1) $existing[$r - 1]
}

# New first entry.
Set-TextValue $boSheet.Cells.Item(1, 1) "0"

# 2. Update the data validation on the RNAseq sheet's barcode_offset column
#    (O) so it covers the new 6-row range instead of the old 5-row range.
$mainSheet = $wb.Worksheets.Item("RNAseq")
$dvRange = $mainSheet.Range("O2:O1001")
$dvRange.Validation.Delete()
$dvRange.Validation.Add(3, 1, 1, "'barcode_offset'!$A$1:$A$6", "")
$dvRange.Validation.ErrorTitle = "Validation Error"
$dvRange.Validation.ErrorMessage = ""
$dvRange.Validation.ShowError = $true

# 3. Bump the schema's pav:createdOn timestamp on the .metadata sheet.
$metaSheet = $wb.Worksheets.Item(".metadata")
$metaSheet.Range("C2").Value = "2023-10-31T14:33:40-07:00"
